# echoMRI editions in timedura column and CD1 exploration updates
#
# 1) Reformat the TimeDateDura (column G) values for rows 2-26 from
#    "HH:MM:SS; D Mon YYYY; NN; ems" to "HH:MM:SS Mon D, YYYY; NN; ems".
# 2) Cosmetic worksheet/view updates: zoom to 170%, select D1, and
#    widen column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = "09:53:53 Aug 4, 2025; 31; ems"
$ws.Range("G3").Value  = "09:55:24 Aug 4, 2025; 32; ems"
$ws.Range("G4").Value  = "09:57:15 Aug 4, 2025; 31; ems"
$ws.Range("G5").Value  = "09:58:56 Aug 4, 2025; 32; ems"
$ws.Range("G6").Value  = "10:00:24 Aug 4, 2025; 27; ems"
$ws.Range("G7").Value  = "10:01:42 Aug 4, 2025; 31; ems"
$ws.Range("G8").Value  = "10:03:14 Aug 4, 2025; 27; ems"
$ws.Range("G9").Value  = "10:05:05 Aug 4, 2025; 32; ems"
$ws.Range("G10").Value = "10:07:13 Aug 4, 2025; 32; ems"
$ws.Range("G11").Value = "10:08:42 Aug 4, 2025; 32; ems"
$ws.Range("G12").Value = "10:10:17 Aug 4, 2025; 41; ems"
$ws.Range("G13").Value = "10:11:50 Aug 4, 2025; 31; ems"
$ws.Range("G14").Value = "10:13:12 Aug 4, 2025; 32; ems"
$ws.Range("G15").Value = "10:14:28 Aug 4, 2025; 31; ems"
$ws.Range("G16").Value = "10:15:43 Aug 4, 2025; 31; ems"
$ws.Range("G17").Value = "10:17:00 Aug 4, 2025; 31; ems"
$ws.Range("G18").Value = "10:18:15 Aug 4, 2025; 31; ems"
$ws.Range("G19").Value = "10:19:38 Aug 4, 2025; 31; ems"
$ws.Range("G20").Value = "10:20:53 Aug 4, 2025; 31; ems"
$ws.Range("G21").Value = "10:22:12 Aug 4, 2025; 31; ems"
$ws.Range("G22").Value = "10:23:28 Aug 4, 2025; 31; ems"
$ws.Range("G23").Value = "10:24:42 Aug 4, 2025; 31; ems"
$ws.Range("G24").Value = "10:25:52 Aug 4, 2025; 27; ems"
$ws.Range("G25").Value = "10:27:01 Aug 4, 2025; 32; ems"
$ws.Range("G26").Value = "10:28:16 Aug 4, 2025; 31; ems"

# Widen column G (TimeDateDura) so the longer reformatted text fits.
$ws.Columns.Item(7).ColumnWidth = 24.6

# Zoom the view to 170% and move the selection to D1.
$excel.ActiveWindow.Zoom = 170
$ws.Range("D1").Select()
